$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.022.75'
$ws.Range("E2").Value = '  -3.70%  '
$ws.Range("D3").Value = '3.227.47'
$ws.Range("E3").Value = '  -4.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.11'
$ws.Range("D5").NumberFormat = "general"
$ws.Range("E5").Value = '  -5.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.48'
$ws.Range("D6").NumberFormat = "general"
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.226.63'
$ws.Range("E8").Value = '  -4.42%  '
$ws.Range("E9").Value = '  -4.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.59'
$ws.Range("D10").NumberFormat = "general"
$ws.Range("E10").Value = '  -5.29%  '
$ws.Range("E11").Value = '  -5.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.397'
$ws.Range("D12").NumberFormat = "general"
$ws.Range("E12").Value = '  -5.14%  '
$ws.Range("D13").Value = '3.779.39'
$ws.Range("E13").Value = '  -4.60%  '
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.16'
$ws.Range("D15").NumberFormat = "general"
$ws.Range("E15").Value = '  -7.94%  '
$ws.Range("D16").Value = '3.224.79'
$ws.Range("E16").Value = '  -4.79%  '
$ws.Range("E17").Value = '  -6.70%  '
$ws.Range("D18").Value = '59.121.79'
$ws.Range("E18").Value = '  -3.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.92'
$ws.Range("D19").NumberFormat = "general"
$ws.Range("E19").Value = '  -7.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.23'
$ws.Range("D20").NumberFormat = "general"
$ws.Range("E20").Value = '  -7.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.32'
$ws.Range("D21").NumberFormat = "general"
$ws.Range("E21").Value = '  -6.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '362.25'
$ws.Range("D22").NumberFormat = "general"
$ws.Range("E22").Value = '  -3.58%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.62'
$ws.Range("D24").NumberFormat = "general"
$ws.Range("E24").Value = '  -6.32%  '
$ws.Range("E25").Value = '  -8.25%  '
$ws.Range("D26").Value = '3.356.68'
$ws.Range("E26").Value = '  -5.26%  '
$ws.Range("D27").Value = '0.0₃0981'
$ws.Range("E27").Value = '  -10.30%  '
$ws.Range("E28").Value = '  -2.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").NumberFormat = "general"
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.06'
$ws.Range("D30").NumberFormat = "general"
$ws.Range("E30").Value = '  -5.45%  '
$ws.Range("E32").Value = '  -7.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.10'
$ws.Range("D33").NumberFormat = "general"
$ws.Range("E33").Value = '  -7.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.99'
$ws.Range("D34").NumberFormat = "general"
$ws.Range("E34").Value = '  -4.46%  '
$ws.Range("E35").Value = '  -2.77%  '
$ws.Range("E36").Value = '  -8.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.48'
$ws.Range("D37").NumberFormat = "general"
$ws.Range("E37").Value = '  -4.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.39'
$ws.Range("D38").NumberFormat = "general"
$ws.Range("E38").Value = '  -6.30%  '
$ws.Range("E39").Value = '  -7.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.91'
$ws.Range("D40").NumberFormat = "general"
$ws.Range("E40").Value = '  -14.95%  '
$ws.Range("E41").Value = '  -7.83%  '
$ws.Range("D42").Value = '3.257.12'
$ws.Range("E42").Value = '  -4.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.06'
$ws.Range("D43").NumberFormat = "general"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.722'
$ws.Range("D44").NumberFormat = "general"
$ws.Range("E44").Value = '  -6.16%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.02'
$ws.Range("D45").NumberFormat = "general"
$ws.Range("E45").Value = '  -7.78%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.09'
$ws.Range("D46").NumberFormat = "general"
$ws.Range("E46").Value = '  -3.89%  '
$ws.Range("E47").Value = '  -6.75%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").Value = '2.290.23'
$ws.Range("E49").Value = '  -9.16%  '
$ws.Range("E50").Value = '  -7.14%  '
$ws.Range("E51").Value = '  -9.68%  '
